$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("DM-Components")
$ws3 = $wb.Worksheets.Item("AE-Components")
$ws4 = $wb.Worksheets.Item("XX-Components")
$ws5 = $wb.Worksheets.Item("CubePrefixes")

# ---------------------------------------------------------------------------
# DM-Components (sheet2): insert two new columns (codeType / nciDomainValue)
# after compName, and populate the NCI codelist info for the SDTM-coded
# components.
# ---------------------------------------------------------------------------
$null = $ws2.Columns("C:D").Insert()
$ws2.Columns("C:D").ColumnWidth = 17.43

$ws2.Range("C1").Value = "codeType"
$ws2.Range("C2").Value = "DATA"
$ws2.Range("C3").Value = "SDTM"
$ws2.Range("C4").Value = "DATA"
$ws2.Range("C5").Value = "DATA"
$ws2.Range("C6").Value = "DATA"
$ws2.Range("C7").Value = "SDTM"

$ws2.Range("D1").Value = "nciDomainValue"
$ws2.Range("D3").Value = "C66731"
$ws2.Range("D7").Value = "C74457"

$null = $ws2.Activate()
$null = $ws2.Range("D12").Select()

# ---------------------------------------------------------------------------
# AE-Components (sheet3): same structural column insert, header labels only.
# ---------------------------------------------------------------------------
$null = $ws3.Columns("C:D").Insert()
$ws3.Columns("C:D").ColumnWidth = 17.43

$ws3.Range("C2").Value = "codeType"
$ws3.Range("D2").Value = "nciDomainValue"

$null = $ws5.Range("A1").Copy()
$null = $ws3.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$null = $ws3.Activate()
$null = $ws3.Range("C2:D2").Select()

# ---------------------------------------------------------------------------
# XX-Components (sheet4): same structural column insert, header labels only.
# ---------------------------------------------------------------------------
$null = $ws4.Columns("C:D").Insert()
$ws4.Columns("C:D").ColumnWidth = 17.43

$ws4.Range("C2").Value = "codeType"
$ws4.Range("D2").Value = "nciDomainValue"

$null = $ws5.Range("A1").Copy()
$null = $ws4.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$null = $ws4.Activate()
$null = $ws4.Range("B41").Select()

# ---------------------------------------------------------------------------
# Intro (sheet1) keeps its own selection, but is no longer the active tab --
# make DM-Components the active sheet/tab last so bookViews/activeTab and
# each sheetView's tabSelected flag land on the right sheet.
# ---------------------------------------------------------------------------
$null = $ws2.Activate()
